$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-525 is a date serial that gets bumped
# by one day (46074 -> 46075) on each automatic update.
$ws.Range("C2:C525").Value = 46075
